$wb = $excel.ActiveWorkbook

$wsWeekly     = $wb.Worksheets.Item("Weekly Task Report")
$wsChart      = $wb.Worksheets.Item("Team Task Chart")
$wsAttendance = $wb.Worksheets.Item("Attendance Records")

# ----------------------------------------------------------------------
# Attendance Records sheet: Week 12 (column M) attendance entries
# ----------------------------------------------------------------------
$wsAttendance.Range("M2").Value = "P"
$wsAttendance.Range("M3").Value = "NP"
$wsAttendance.Range("M4").Value = "P"
$wsAttendance.Range("M5").Value = "NP"
$wsAttendance.Range("M6").Value = "P"

# ----------------------------------------------------------------------
# Team Task Chart sheet
# ----------------------------------------------------------------------
# Updated progress percentages
$wsChart.Range("H9").Value = 1
$wsChart.Range("H12").Value = 0.95
$wsChart.Range("H13").Value = 0.8

# Mark Week 12 (columns BL:BP) complete with "X" for tasks in rows 9-13
$wsChart.Range("BL9:BP9").Value = "X"
$wsChart.Range("BL10:BP10").Value = "X"
$wsChart.Range("BL11:BP11").Value = "X"
$wsChart.Range("BL12:BP12").Value = "X"
$wsChart.Range("BL13:BP13").Value = "X"

# ----------------------------------------------------------------------
# Weekly Task Report sheet
# ----------------------------------------------------------------------
# Week Number: 11 -> 12
$wsWeekly.Range("B7").Value = 12

# Overall team status: At Risk -> On Track
$wsWeekly.Range("B12").Value = "On Track"

# Status narrative text
$wsWeekly.Range("A14").Value = "We have milestone 1 due Sunday (drive system). I've been working on the report. The team has gotten the control system working (getting the control signal from source to destination and translating between forms). We don't have a chassis, suspension, or E-box, so we'll be getting a bit creative with the video to show it working."

# Status labels that mirror the Team Task Chart rows (manually maintained, not formulas)
$wsWeekly.Range("D30").Value = "On Track"
$wsWeekly.Range("D31").Value = "On Track"

# ----------------------------------------------------------------------
# Restore sheet view selections (cursor moved since last save)
# ----------------------------------------------------------------------
$wsChart.Activate()
$wsChart.Range("H10").Select()

$wsAttendance.Activate()
$wsAttendance.Range("M6").Select()

$wsWeekly.Activate()
$wsWeekly.Range("A21").Select()
